# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (new week: 2022-06-20 / serial 44740) above
# the current row 63, pushing the existing data block (old rows 63:162) down
# to rows 66:165, then populate the 3 newly inserted rows with the new
# Especial / Primera / Segunda quality triplet for Kiwi Hayward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 63; this shifts old rows 63-162 down to
# become rows 66-165 (dimension grows from T162 to T165 automatically).
$ws.Rows("63:65").Insert()

# Shared context values (identical across the three new rows).
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad = "Hayward"
$unidad = "`$/bandeja 18 kilos"
$origen = "Región de O'Higgins"
$kgUnidad = 18
$fecha = 44740

function Set-KiwiRow($row, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-KiwiRow 63 "Especial" 50 10000 10000 10000 556
Set-KiwiRow 64 "Primera"  50 8000  8000  8000  444
Set-KiwiRow 65 "Segunda"  50 6000  6000  6000  333
